$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the existing ListObject (Table1) to include two new columns (D:E)
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:E11"))

# Name the new header cells through the table's header row so the table
# definition (tableColumn name=...) picks up the new names, not just the
# literal cell text.
$table.HeaderRowRange.Cells.Item(1, 4).Value = "ISSELECTED"
$table.HeaderRowRange.Cells.Item(1, 5).Value = "DOORPRICENAME"

# Update the active selection to match the target state
$ws.Range("E2").Select() | Out-Null
